$d = $word.ActiveDocument

$null = $d.Content.Find.Execute("2023-01-28 Saturday", $true, $true, $false, $false, $false, $true, 1, $false, "2023-01-29 Sunday", 2)
$null = $d.Content.Find.Execute("19+65=", $true, $true, $false, $false, $false, $true, 1, $false, "38+42=", 2)
$null = $d.Content.Find.Execute("32-14=", $true, $true, $false, $false, $false, $true, 1, $false, "35+17=", 2)
$null = $d.Content.Find.Execute("9+16=", $true, $true, $false, $false, $false, $true, 1, $false, "69-28=", 2)
$null = $d.Content.Find.Execute("79+17=", $true, $true, $false, $false, $false, $true, 1, $false, "60-23=", 2)
$null = $d.Content.Find.Execute("85+9=", $true, $true, $false, $false, $false, $true, 1, $false, "94-87=", 2)
$null = $d.Content.Find.Execute("9+50=", $true, $true, $false, $false, $false, $true, 1, $false, "28+30=", 2)
$null = $d.Content.Find.Execute("89-5=", $true, $true, $false, $false, $false, $true, 1, $false, "65-62=", 2)
$null = $d.Content.Find.Execute("55+10=", $true, $true, $false, $false, $false, $true, 1, $false, "47-24=", 2)
$null = $d.Content.Find.Execute("72-49=", $true, $true, $false, $false, $false, $true, 1, $false, "99-27=", 2)
$null = $d.Content.Find.Execute("8+41=", $true, $true, $false, $false, $false, $true, 1, $false, "37+35=", 2)
$null = $d.Content.Find.Execute("79-4=", $true, $true, $false, $false, $false, $true, 1, $false, "1+41=", 2)
$null = $d.Content.Find.Execute("3+32=", $true, $true, $false, $false, $false, $true, 1, $false, "15+32=", 2)
$null = $d.Content.Find.Execute("91-13=", $true, $true, $false, $false, $false, $true, 1, $false, "87-23=", 2)
$null = $d.Content.Find.Execute("28-26=", $true, $true, $false, $false, $false, $true, 1, $false, "51+21=", 2)
$null = $d.Content.Find.Execute("50+42=", $true, $true, $false, $false, $false, $true, 1, $false, "0+92=", 2)
$null = $d.Content.Find.Execute("19+63=", $true, $true, $false, $false, $false, $true, 1, $false, "64-43=", 2)
$null = $d.Content.Find.Execute("61+0=", $true, $true, $false, $false, $false, $true, 1, $false, "57-29=", 2)
$null = $d.Content.Find.Execute("0+86=", $true, $true, $false, $false, $false, $true, 1, $false, "95-47=", 2)
$null = $d.Content.Find.Execute("15+71=", $true, $true, $false, $false, $false, $true, 1, $false, "86-39=", 2)
$null = $d.Content.Find.Execute("38-9=", $true, $true, $false, $false, $false, $true, 1, $false, "51-39=", 2)
$null = $d.Content.Find.Execute("19-4=", $true, $true, $false, $false, $false, $true, 1, $false, "27+17=", 2)
$null = $d.Content.Find.Execute("24+6=", $true, $true, $false, $false, $false, $true, 1, $false, "83-45=", 2)
$null = $d.Content.Find.Execute("36-19=", $true, $true, $false, $false, $false, $true, 1, $false, "3+73=", 2)
$null = $d.Content.Find.Execute("74-32=", $true, $true, $false, $false, $false, $true, 1, $false, "51+28=", 2)
$null = $d.Content.Find.Execute("63-0=", $true, $true, $false, $false, $false, $true, 1, $false, "53-24=", 2)
$null = $d.Content.Find.Execute("50+15=", $true, $true, $false, $false, $false, $true, 1, $false, "57-36=", 2)
$null = $d.Content.Find.Execute("42-8=", $true, $true, $false, $false, $false, $true, 1, $false, "74+23=", 2)
$null = $d.Content.Find.Execute("66-2=", $true, $true, $false, $false, $false, $true, 1, $false, "67-57=", 2)
$null = $d.Content.Find.Execute("22+54=", $true, $true, $false, $false, $false, $true, 1, $false, "93-91=", 2)
$null = $d.Content.Find.Execute("81+15=", $true, $true, $false, $false, $false, $true, 1, $false, "9+59=", 2)
$null = $d.Content.Find.Execute("75-51=", $true, $true, $false, $false, $false, $true, 1, $false, "88-29=", 2)
$null = $d.Content.Find.Execute("50-37=", $true, $true, $false, $false, $false, $true, 1, $false, "2+13=", 2)
$null = $d.Content.Find.Execute("46+45=", $true, $true, $false, $false, $false, $true, 1, $false, "28+3=", 2)
$null = $d.Content.Find.Execute("46+52=", $true, $true, $false, $false, $false, $true, 1, $false, "94-14=", 2)
$null = $d.Content.Find.Execute("75-61=", $true, $true, $false, $false, $false, $true, 1, $false, "34+5=", 2)
$null = $d.Content.Find.Execute("15+78=", $true, $true, $false, $false, $false, $true, 1, $false, "87-13=", 2)
$null = $d.Content.Find.Execute("53+12=", $true, $true, $false, $false, $false, $true, 1, $false, "81-47=", 2)
$null = $d.Content.Find.Execute("69-10=", $true, $true, $false, $false, $false, $true, 1, $false, "77-49=", 2)
$null = $d.Content.Find.Execute("14+16=", $true, $true, $false, $false, $false, $true, 1, $false, "96-14=", 2)
$null = $d.Content.Find.Execute("55-49=", $true, $true, $false, $false, $false, $true, 1, $false, "79+5=", 2)
$null = $d.Content.Find.Execute("59-56=", $true, $true, $false, $false, $false, $true, 1, $false, "7+40=", 2)
$null = $d.Content.Find.Execute("94-62=", $true, $true, $false, $false, $false, $true, 1, $false, "49-8=", 2)
$null = $d.Content.Find.Execute("64-1=", $true, $true, $false, $false, $false, $true, 1, $false, "49+22=", 2)
$null = $d.Content.Find.Execute("51+31=", $true, $true, $false, $false, $false, $true, 1, $false, "48+45=", 2)
$null = $d.Content.Find.Execute("8+28=", $true, $true, $false, $false, $false, $true, 1, $false, "25-22=", 2)
$null = $d.Content.Find.Execute("51+11=", $true, $true, $false, $false, $false, $true, 1, $false, "50-21=", 2)
$null = $d.Content.Find.Execute("98-75=", $true, $true, $false, $false, $false, $true, 1, $false, "2+10=", 2)
$null = $d.Content.Find.Execute("57+38=", $true, $true, $false, $false, $false, $true, 1, $false, "22+76=", 2)
$null = $d.Content.Find.Execute("0+95=", $true, $true, $false, $false, $false, $true, 1, $false, "23-12=", 2)
$null = $d.Content.Find.Execute("45+29=", $true, $true, $false, $false, $false, $true, 1, $false, "70-4=", 2)
$null = $d.Content.Find.Execute("67-55=", $true, $true, $false, $false, $false, $true, 1, $false, "93-30=", 2)
$null = $d.Content.Find.Execute("86-30=", $true, $true, $false, $false, $false, $true, 1, $false, "84-69=", 2)
$null = $d.Content.Find.Execute("77-2=", $true, $true, $false, $false, $false, $true, 1, $false, "59+28=", 2)
$null = $d.Content.Find.Execute("96-7=", $true, $true, $false, $false, $false, $true, 1, $false, "54+18=", 2)
$null = $d.Content.Find.Execute("82-3=", $true, $true, $false, $false, $false, $true, 1, $false, "79+9=", 2)
$null = $d.Content.Find.Execute("41+27=", $true, $true, $false, $false, $false, $true, 1, $false, "2+10=", 2)
$null = $d.Content.Find.Execute("74+21=", $true, $true, $false, $false, $false, $true, 1, $false, "90-28=", 2)
$null = $d.Content.Find.Execute("98-36=", $true, $true, $false, $false, $false, $true, 1, $false, "13-7=", 2)
$null = $d.Content.Find.Execute("41+21=", $true, $true, $false, $false, $false, $true, 1, $false, "97-63=", 2)
$null = $d.Content.Find.Execute("92-9=", $true, $true, $false, $false, $false, $true, 1, $false, "4+82=", 2)
$null = $d.Content.Find.Execute("20+32=", $true, $true, $false, $false, $false, $true, 1, $false, "85-0=", 2)
$null = $d.Content.Find.Execute("98-0=", $true, $true, $false, $false, $false, $true, 1, $false, "46+33=", 2)
$null = $d.Content.Find.Execute("66-37=", $true, $true, $false, $false, $false, $true, 1, $false, "8+70=", 2)
$null = $d.Content.Find.Execute("31-18=", $true, $true, $false, $false, $false, $true, 1, $false, "35+29=", 2)
$null = $d.Content.Find.Execute("28+28=", $true, $true, $false, $false, $false, $true, 1, $false, "81+9=", 2)
$null = $d.Content.Find.Execute("82-40=", $true, $true, $false, $false, $false, $true, 1, $false, "1+31=", 2)
$null = $d.Content.Find.Execute("27-6=", $true, $true, $false, $false, $false, $true, 1, $false, "36-29=", 2)
$null = $d.Content.Find.Execute("99-64=", $true, $true, $false, $false, $false, $true, 1, $false, "49+17=", 2)
$null = $d.Content.Find.Execute("0+19=", $true, $true, $false, $false, $false, $true, 1, $false, "94-45=", 2)
$null = $d.Content.Find.Execute("65+3=", $true, $true, $false, $false, $false, $true, 1, $false, "48+29=", 2)
$null = $d.Content.Find.Execute("90-24=", $true, $true, $false, $false, $false, $true, 1, $false, "80-26=", 2)
$null = $d.Content.Find.Execute("8-6=", $true, $true, $false, $false, $false, $true, 1, $false, "36+4=", 2)
$null = $d.Content.Find.Execute("91+7=", $true, $true, $false, $false, $false, $true, 1, $false, "50+26=", 2)
$null = $d.Content.Find.Execute("38-21=", $true, $true, $false, $false, $false, $true, 1, $false, "98-50=", 2)
$null = $d.Content.Find.Execute("88-27=", $true, $true, $false, $false, $false, $true, 1, $false, "17+46=", 2)
$null = $d.Content.Find.Execute("26-24=", $true, $true, $false, $false, $false, $true, 1, $false, "41+30=", 2)
$null = $d.Content.Find.Execute("16+23=", $true, $true, $false, $false, $false, $true, 1, $false, "17+13=", 2)
$null = $d.Content.Find.Execute("16+2=", $true, $true, $false, $false, $false, $true, 1, $false, "33+11=", 2)
$null = $d.Content.Find.Execute("49+38=", $true, $true, $false, $false, $false, $true, 1, $false, "62-50=", 2)
$null = $d.Content.Find.Execute("30+38=", $true, $true, $false, $false, $false, $true, 1, $false, "57+16=", 2)
$null = $d.Content.Find.Execute("6+7=", $true, $true, $false, $false, $false, $true, 1, $false, "81-51=", 2)
$null = $d.Content.Find.Execute("78-68=", $true, $true, $false, $false, $false, $true, 1, $false, "64-61=", 2)
$null = $d.Content.Find.Execute("10+51=", $true, $true, $false, $false, $false, $true, 1, $false, "22+50=", 2)
$null = $d.Content.Find.Execute("95-78=", $true, $true, $false, $false, $false, $true, 1, $false, "63+25=", 2)
$null = $d.Content.Find.Execute("39-24=", $true, $true, $false, $false, $false, $true, 1, $false, "87-65=", 2)
$null = $d.Content.Find.Execute("78+2=", $true, $true, $false, $false, $false, $true, 1, $false, "17+0=", 2)
$null = $d.Content.Find.Execute("10+12=", $true, $true, $false, $false, $false, $true, 1, $false, "19+68=", 2)
$null = $d.Content.Find.Execute("63+5=", $true, $true, $false, $false, $false, $true, 1, $false, "96-8=", 2)
$null = $d.Content.Find.Execute("24+65=", $true, $true, $false, $false, $false, $true, 1, $false, "57+41=", 2)
$null = $d.Content.Find.Execute("35+9=", $true, $true, $false, $false, $false, $true, 1, $false, "55+13=", 2)
$null = $d.Content.Find.Execute("11-8=", $true, $true, $false, $false, $false, $true, 1, $false, "34-23=", 2)
$null = $d.Content.Find.Execute("11+36=", $true, $true, $false, $false, $false, $true, 1, $false, "30+53=", 2)
$null = $d.Content.Find.Execute("91-37=", $true, $true, $false, $false, $false, $true, 1, $false, "69-49=", 2)
$null = $d.Content.Find.Execute("15+68=", $true, $true, $false, $false, $false, $true, 1, $false, "61-38=", 2)
$null = $d.Content.Find.Execute("90-78=", $true, $true, $false, $false, $false, $true, 1, $false, "13+62=", 2)
$null = $d.Content.Find.Execute("95-23=", $true, $true, $false, $false, $false, $true, 1, $false, "9+17=", 2)
$null = $d.Content.Find.Execute("84-71=", $true, $true, $false, $false, $false, $true, 1, $false, "61-27=", 2)
$null = $d.Content.Find.Execute("91-81=", $true, $true, $false, $false, $false, $true, 1, $false, "4+48=", 2)
$null = $d.Content.Find.Execute("10+88=", $true, $true, $false, $false, $false, $true, 1, $false, "33-6=", 2)
$null = $d.Content.Find.Execute("14+5=", $true, $true, $false, $false, $false, $true, 1, $false, "12+82=", 2)
